$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that is bumped by one
# day (45204 -> 45205, i.e. 2023-10-05 -> 2023-10-06) for every data row,
# from row 2 through the last used row (518).
$ws.Range("C2:C518").Value = 45205
